# The sheet contains a table of price observations for "Feria Lagunitas de
# Puerto Montt - Cilantro". A new weekly observation is inserted as a new
# row 276 (pushing the existing rows 276-290 down to 277-291), and the
# sheet dimension grows from A1:R290 to A1:R291.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row above the current row 276 - this shifts the old
# rows 276..290 down to 277..291 automatically.
$ws.Rows(276).Insert()

# Populate the newly inserted row 276 with the new observation.
$ws.Range("A276").Value = 4
$ws.Range("B276").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C276").Value = "Los Lagos"
$ws.Range("D276").Value = 44706
$ws.Range("E276").Value = 10
$ws.Range("F276").Value = 100112040
$ws.Range("G276").Value = "Cilantro"
$ws.Range("H276").Value = "Sin especificar"
$ws.Range("I276").Value = "Primera"
$ws.Range("J276").Value = 25
$ws.Range("K276").Value = 6000
$ws.Range("L276").Value = 6000
$ws.Range("M276").Value = 6000
$ws.Range("N276").Value = "$/docena de atados (2 kilos)"
$ws.Range("O276").Value = "Región de La Araucanía"
$ws.Range("P276").Value = 3000
$ws.Range("Q276").Value = 2
$ws.Range("R276").Value = "Hortaliza"
